# Cập nhật danh sách chức năng
# Applies the data + view changes recorded in the commit to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate() | Out-Null

# --- Cell value updates (column E "percent done" + G15 name swap) ---

# G15: reassign the responsible person from "Huy" to "Nhi"
$ws.Range("G15").Value = "Nhi"

# Column E (progress) newly filled in for several rows
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("E47").Value = 0.3

# --- View state: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("G14").Select() | Out-Null

Write-Output "Applied function-list update"
